$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at position 836, pushing existing rows 836+ down to 838+.
$ws.Range("A836:A837").EntireRow.Insert()

# New weekly record row 836
$ws.Range("A836").Value = 5
$ws.Range("B836").Value = "Macroferia Regional de Talca"
$ws.Range("C836").Value = "Maule"
$ws.Range("D836").Value2 = 44984
$ws.Range("E836").Value = 7
$ws.Range("F836").Value = "Fruta"
$ws.Range("G836").Value = 100108
$ws.Range("H836").Value = "Tropicales y subtropicales"
$ws.Range("I836").Value = 100108006
$ws.Range("J836").Value = "Plátano"
$ws.Range("K836").Value = "Sin especificar"
$ws.Range("L836").Value = "Pintón"
$ws.Range("M836").Value = 1000
$ws.Range("N836").Value = 22000
$ws.Range("O836").Value = 22000
$ws.Range("P836").Value = 22000
$ws.Range("Q836").Value = "$/caja 20 kilos"
$ws.Range("R836").Value = "Ecuador"
$ws.Range("S836").Value = 1100
$ws.Range("T836").Value = 20

# New weekly record row 837
$ws.Range("A837").Value = 5
$ws.Range("B837").Value = "Macroferia Regional de Talca"
$ws.Range("C837").Value = "Maule"
$ws.Range("D837").Value2 = 44984
$ws.Range("E837").Value = 7
$ws.Range("F837").Value = "Fruta"
$ws.Range("G837").Value = 100108
$ws.Range("H837").Value = "Tropicales y subtropicales"
$ws.Range("I837").Value = 100108006
$ws.Range("J837").Value = "Plátano"
$ws.Range("K837").Value = "Sin especificar"
$ws.Range("L837").Value = "Primera Pintón"
$ws.Range("M837").Value = 700
$ws.Range("N837").Value = 23000
$ws.Range("O837").Value = 24000
$ws.Range("P837").Value = 23429
$ws.Range("Q837").Value = "$/caja 20 kilos"
$ws.Range("R837").Value = "Ecuador"
$ws.Range("S837").Value = 1171
$ws.Range("T837").Value = 20

# Ensure date formatting on the new D cells matches the rest of the column.
$ws.Range("D836:D837").NumberFormat = "YYYY-MM-DD HH:MM:SS"
